$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.098.38"
$ws.Range("E2").Value = "  -0.69%  "

$ws.Range("D3").Value = "3.291.28"
$ws.Range("E3").Value = "  -0.79%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'585.36"
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").Value = "'181.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.56%  "

$ws.Range("E7").Value = "  +8.27%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -3.09%  "

$ws.Range("E10").Value = "  +1.64%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "3.862.64"
$ws.Range("E12").Value = "  -0.83%  "

$ws.Range("E13").Value = "  -4.58%  "

$ws.Range("D14").Value = "66.124.49"
$ws.Range("E14").Value = "  -0.75%  "

$ws.Range("D15").Value = "'26.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.47%  "

$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").Value = "3.264.75"
$ws.Range("E17").Value = "  -1.83%  "

$ws.Range("D18").Value = "'433.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.35%  "

$ws.Range("D19").Value = "'13.27"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.26%  "

$ws.Range("D20").Value = "'5.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.21%  "

$ws.Range("E21").Value = "  -3.07%  "

$ws.Range("D22").Value = "'72.32"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.99%  "

$ws.Range("E23").Value = "  +0.05%  "

$ws.Range("D24").Value = "'5.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").Value = "3.433.69"
$ws.Range("E25").Value = "  -0.74%  "

$ws.Range("D26").Value = "'0.511"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.79%  "

$ws.Range("D27").Value = "'0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.04%  "

$ws.Range("E28").Value = "  +2.10%  "

$ws.Range("D29").Value = "'8.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.82%  "

$ws.Range("D30").Value = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "

$ws.Range("E31").Value = "  +0.63%  "

$ws.Range("D32").Value = "'22.35"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.36%  "

$ws.Range("E33").Value = "  -0.03%  "

$ws.Range("D34").Value = "'5.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.82%  "

$ws.Range("D35").Value = "'6.62"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.22%  "

$ws.Range("E36").Value = "  -2.52%  "

$ws.Range("D37").Value = "'158.51"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.13%  "

$ws.Range("E38").Value = "  -5.25%  "

$ws.Range("D39").Value = "'26.51"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.85%  "

$ws.Range("E40").Value = "  -3.57%  "

$ws.Range("D41").Value = "2.780.75"
$ws.Range("E41").Value = "  -1.05%  "

$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("D43").Value = "'4.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.74%  "

$ws.Range("D44").Value = "'40.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").Value = "'6.02"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.29%  "

$ws.Range("D46").Value = "'0.0660"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.49%  "

$ws.Range("D47").Value = "'2.32"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.16%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'23.29"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.92%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").Value = "'316.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.69%  "

$ws.Range("E50").Value = "  -1.61%  "

$ws.Range("E51").Value = "  +5.75%  "
